# Add a new "Slovakia" worksheet based on the existing "Portugal" sheet,
# fill in the Slovakia-specific values, and update sheet selections/active tab.

$wb = $excel.ActiveWorkbook

$portugal = $wb.Worksheets.Item("Portugal")

# Duplicate the Portugal sheet and place the copy right after it.
[void]$portugal.Copy($null, $portugal)
$slovakia = $wb.Worksheets.Item($wb.Worksheets.Count)
$slovakia.Name = "Slovakia"

# Update the market name and Jira/NGC reference for Slovakia.
$slovakia.Range("B2").Value = "Slovakia Market"
$slovakia.Range("B4").Value = "NGC-2930/T3177/T3176/T3179/T3178"

# Reset the row heights on the new sheet back to the default (Portugal's
# rows 3-5 were manually resized; the new sheet should not inherit that).
[void]$slovakia.Rows("3:5").AutoFit()

# The Portugal sheet is no longer the active tab; its selection becomes
# the whole sheet.
[void]$portugal.Cells.Select()

# Make Slovakia the active sheet/tab with B9 selected.
[void]$slovakia.Select()
[void]$slovakia.Range("B9").Select()
